# Adds ~30 new Spanish/Kawaiinese vocabulary rows (122-151) to the
# dictionary sheet, mirroring the "Add files via upload" commit that
# appended new <si> shared strings + sheetData rows, plus one new font
# (dark-gray, rgb FF0D0D0D) applied only to the Kawaiinese translation
# of "Cuantos" (row 136, column B).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: row number, Spanish (col A), Kawaiinese (col B), optional
# attribution note for col C ("ChatGPT" on several of the newer rows).
$rows = @(
    @(122, 'Encender', 'onu', ''),
    @(123, 'Rapido', 'Fatsu', ''),
    @(124, 'Rapidamentre', 'Fatsumenu', ''),
    @(125, 'Mucho', 'Arotu', ''),
    @(126, 'Tambien', 'Chu', ''),
    @(127, 'Descansar', 'Driamu', ''),
    @(128, 'Descanso', 'Driamo', ''),
    @(129, 'Libro', 'Honu', ''),
    @(130, 'Este', 'Kono', 'ChatGPT'),
    @(131, 'Ese', 'Ano', 'ChatGPT'),
    @(132, 'Aquel', 'Sono', 'ChatGPT'),
    @(133, 'Esto', 'Kore', 'ChatGPT'),
    @(134, 'Eso', 'Sore', 'ChatGPT'),
    @(135, 'Aquello', 'Are', 'ChatGPT'),
    @(136, 'Cuantos', 'Ikutsu', 'ChatGPT'),
    @(137, 'Alguien', 'Dareka', 'ChatGPT'),
    @(138, 'Algo', 'Nanika', 'ChatGPT'),
    @(139, 'Igual', 'Iruha', ''),
    @(140, 'Mas que', 'Miruha', ''),
    @(141, 'El mas', 'Supra', ''),
    @(142, 'Uno mismo', 'Jibun', 'ChatGPT'),
    @(143, 'Si mismo', 'Jishin', 'ChatGPT'),
    @(144, 'Yo misma', 'Atashi', 'ChatGPT'),
    @(145, 'De uno mismo', 'Jibun no', 'ChatGPT'),
    @(146, 'De si mismo', 'Jishin no', 'ChatGPT'),
    @(147, 'Si (comp.)', 'Ifu', ''),
    @(148, 'Entonces', 'Sou', ''),
    @(149, 'Aprobar', 'Ganaru', ''),
    @(150, 'Musico', 'Takimusianu', ''),
    @(151, 'De', 'Oku', '')
)

foreach ($entry in $rows) {
    $r = $entry[0]
    $ws.Cells.Item($r, 1).Value = $entry[1]
    $ws.Cells.Item($r, 2).Value = $entry[2]
    if ($entry[3] -ne '') {
        $ws.Cells.Item($r, 3).Value = $entry[3]
    }
}

# Row 136's Kawaiinese word ("Ikutsu") gets its own dark-gray font
# (new 3rd font entry in styles.xml), distinguishing it from the rest.
$ws.Range("B136").Font.Color = 855309

# Reflect where the editor had scrolled/clicked last.
$ws.Range("A46").Select()
